{"js": "// The commit removes 5 consecutive paragraphs that immediately follow the\n// paragraph containing the text \"GDWFSDVSWEBAJAVAEXAIII1A\":\n//   1. an empty paragraph (centered)\n//   2. an empty paragraph (shaded, left aligned)\n//   3. the paragraph \"Ceci est un mod\u00e8le de copie. N'oubliez pas de\n//      renseigner vos pr\u00e9nom/nom, ainsi que le nom et le lien vers le\n//      projet.\"\n//   4. an empty paragraph (shaded, left aligned)\n//   5. the paragraph \"Vous pouvez bien s\u00fbr agrandir les cadres pour\n//      r\u00e9pondre aux questions sur la description du projet si\n//      n\u00e9cessaire.\"\n//\n// Locate the anchor paragraph by its (unique) text, then delete the five\n// paragraphs that come right after it.\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"GDWFSDVSWEBAJAVAEXAIII1A\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\n\nconst paragraphsToRemoveCount = 5;\nlet current = anchorParagraph.getNext();\nfor (let i = 0; i < paragraphsToRemoveCount; i++) {\n  const toDelete = current;\n  // Grab the following paragraph before deleting the current one so the\n  // chain stays valid once `toDelete` is gone.\n  current = toDelete.getNext();\n  toDelete.delete();\n}\n\nawait context.sync();\n", "ps1": "# The commit removes 5 consecutive paragraphs that immediately follow the\n# paragraph containing the text \"GDWFSDVSWEBAJAVAEXAIII1A\":\n#   1. an empty paragraph (centered)\n#   2. an empty paragraph (shaded, left aligned)\n#   3. the paragraph \"Ceci est un modele de copie. N'oubliez pas de\n#      renseigner vos prenom/nom, ainsi que le nom et le lien vers le\n#      projet.\"\n#   4. an empty paragraph (shaded, left aligned)\n#   5. the paragraph \"Vous pouvez bien sur agrandir les cadres pour\n#      repondre aux questions sur la description du projet si\n#      necessaire.\"\n#\n# Locate the anchor paragraph by its (unique) text, then delete the five\n# paragraphs that come right after it.\n\n$d = $word.ActiveDocument\n\n$anchorRange = $d.Content\n$find = $anchorRange.Find\n$find.Text = \"GDWFSDVSWEBAJAVAEXAIII1A\"\n$found = $find.Execute()\n\n# Read the paragraph index off the SAME range object the Find ran\n# against (re-fetching $d.Content here would yield a fresh, un-positioned\n# range starting back at the top of the document).\n$anchorIndex = $anchorRange.Paragraphs.First.Index\n\n$paragraphsToRemoveCount = 5\nfor ($i = 0; $i -lt $paragraphsToRemoveCount; $i++) {\n    $p = $d.Paragraphs.Item($anchorIndex + 1)\n    $p.Range.Delete()\n}\n"}
